$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.051.57'
$ws.Range("E2").Value = '  +2.42%  '

$ws.Range("D3").Value = '1.675.05'
$ws.Range("E3").Value = '  +3.47%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = "'216.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.61%  '

$ws.Range("E6").Value = '  +2.16%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  +2.56%  '

$ws.Range("D9").Value = "'0.0618"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.54%  '

$ws.Range("D10").Value = "'20.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.20%  '

$ws.Range("D11").Value = "'0.0888"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.79%  '

$ws.Range("D12").Value = '1.910.31'
$ws.Range("E12").Value = '  +3.42%  '

$ws.Range("D13").Value = '1.673.17'
$ws.Range("E13").Value = '  +3.33%  '

$ws.Range("E14").Value = '  +1.79%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = "'0.522"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.65%  '

$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = "'65.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.25%  '

$ws.Range("D17").Value = '27.057.88'
$ws.Range("E17").Value = '  +2.37%  '

$ws.Range("D18").Value = "'236.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.02%  '

$ws.Range("E19").Value = '  +1.76%  '

$ws.Range("E20").Value = '  -0.98%  '

$ws.Range("E21").Value = '  +0.08%  '

$ws.Range("D22").Value = "'4.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.09%  '

$ws.Range("D23").Value = "'9.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.88%  '

$ws.Range("E24").Value = '  +1.43%  '

$ws.Range("D25").Value = "'145.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.09%  '

$ws.Range("E26").Value = '  +1.48%  '

$ws.Range("E27").Value = '  +0.50%  '

$ws.Range("E28").Value = '  +2.49%  '

$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("D30").Value = "'0.0499"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.47%  '

$ws.Range("E31").Value = '  +1.96%  '

$ws.Range("E32").Value = '  +1.97%  '

$ws.Range("D33").Value = '1.472.85'
$ws.Range("E33").Value = '  -2.82%  '

$ws.Range("E34").Value = '  +5.23%  '

$ws.Range("E35").Value = '  +6.75%  '

$ws.Range("E36").Value = '  -0.60%  '

$ws.Range("E37").Value = '  +1.03%  '

$ws.Range("D38").Value = "'0.898"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.58%  '

$ws.Range("E39").Value = '  +2.09%  '

$ws.Range("D40").Value = "'6.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.23%  '

$ws.Range("E41").Value = '  +0.07%  '

$ws.Range("E42").Value = '  +10.82%  '

$ws.Range("E43").Value = '  +3.33%  '

$ws.Range("D44").Value = "'66.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.33%  '

$ws.Range("D45").Value = '1.819.23'
$ws.Range("E45").Value = '  +3.37%  '

$ws.Range("E46").Value = '  +2.08%  '

$ws.Range("D47").Value = "'90.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.36%  '

$ws.Range("E48").Value = '  +2.03%  '

$ws.Range("E49").Value = '  -0.94%  '

$ws.Range("E50").Value = '  +4.04%  '

$ws.Range("E51").Value = '  +1.15%  '
